$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 18:35"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1652076
$ws.Range("C4").Value = 6982
$ws.Range("D4").Value = 404283
$ws.Range("E4").Value = 1149943
$ws.Range("G4").Value = 203
$ws.Range("H4").Value = 97850

# Italia (row 9)
$ws.Range("B9").Value = 229327
$ws.Range("C9").Value = 669
$ws.Range("D9").Value = 138840
$ws.Range("E9").Value = 57752
$ws.Range("G9").Value = 119
$ws.Range("H9").Value = 32735

# India (row 14)
$ws.Range("B14").Value = 130908
$ws.Range("C14").Value = 6114
$ws.Range("D14").Value = 54179
$ws.Range("E14").Value = 72869
$ws.Range("G14").Value = 134
$ws.Range("H14").Value = 3860

# Irlanda (row 33)
$ws.Range("B33").Value = 24582
$ws.Range("C33").Value = 76
$ws.Range("E33").Value = 1918
$ws.Range("G33").Value = 12
$ws.Range("H33").Value = 1604

# Chequia (row 53)
$ws.Range("B53").Value = 8853
$ws.Range("C53").Value = 40
$ws.Range("D53").Value = 6038
$ws.Range("E53").Value = 2501
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 314

# Marruecos (row 58)
$ws.Range("B58").Value = 7406
$ws.Range("C58").Value = 74
$ws.Range("D58").Value = 4638
$ws.Range("E58").Value = 2570
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 198

# Jordania (row 124) - values updated, country stays
$ws.Range("B124").Value = 704
$ws.Range("C124").Value = 4
$ws.Range("D124").Value = 470
$ws.Range("E124").Value = 225

# Row 125 (San Marino) unchanged

# Rows 126-128 reorder: Sierra Leona moves up above Chad & Malta (driven by its new, higher case count)
# Row 126 becomes Sierra Leona with fresh data
$ws.Range("A126").Value = "Sierra Leona"
$ws.Range("B126").Value = 621
$ws.Range("C126").Value = 15
$ws.Range("D126").Value = 241
$ws.Range("E126").Value = 341
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 39

# Row 127 becomes Republica del Chad (its data is unchanged, just shifted down one row)
$ws.Range("A127").Value = "Republica del Chad"
$ws.Range("B127").Value = 611
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 196
$ws.Range("E127").Value = 357
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 58

# Row 128 becomes Malta (its data is unchanged, just shifted down one row)
$ws.Range("A128").Value = "Malta"
$ws.Range("B128").Value = 609
$ws.Range("C128").Value = 9
$ws.Range("D128").Value = 473
$ws.Range("E128").Value = 130
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 6

# Row 129 (Nepal) unchanged

# Libia (row 177)
$ws.Range("B177").Value = 75
$ws.Range("C177").Value = 3
$ws.Range("D177").Value = 39
$ws.Range("E177").Value = 33

# Curazao (row 202)
$ws.Range("B202").Value = 17
$ws.Range("C202").Value = 1
$ws.Range("E202").Value = 2
